# Auto-generated edit script: updates crypto price/volume columns (D, E)
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = "@"
$cell.Value = '27.293.10'
$cell.Style = "Normal"
$ws.Range('E2').Value = '  -4.02%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = "@"
$cell.Value = '1.845.68'
$cell.Style = "Normal"
$ws.Range('E3').Value = '  -5.67%  '
$ws.Range('E4').Value = '  -0.41%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '321.17'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('E6').Value = '  -0.37%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = "@"
$cell.Value = '0.4456'
$cell.Style = "Normal"
$ws.Range('E7').Value = '  -6.38%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '0.3817'
$cell.Style = "Normal"
$ws.Range('E8').Value = '  -5.80%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '48.56'
$cell.Style = "Normal"
$ws.Range('E9').Value = '  -8.99%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '0.07839'
$cell.Style = "Normal"
$ws.Range('E10').Value = '  -6.75%  '
$ws.Range('E11').Value = '  -4.76%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '21.27'
$cell.Style = "Normal"
$ws.Range('E12').Value = '  -4.65%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '1.839.44'
$cell.Style = "Normal"
$ws.Range('E13').Value = '  -5.59%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '5.833'
$cell.Style = "Normal"
$ws.Range('E14').Value = '  -5.11%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '7.069'
$cell.Style = "Normal"
$ws.Range('E15').Value = '  -7.11%  '
$ws.Range('E16').Value = '  -0.44%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '85.57'
$cell.Style = "Normal"
$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '0.00001023'
$cell.Style = "Normal"
$ws.Range('E18').Value = '  -4.13%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '0.06503'
$cell.Style = "Normal"
$ws.Range('E19').Value = '  -1.34%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '16.92'
$cell.Style = "Normal"
$ws.Range('E20').Value = '  -8.62%  '
$ws.Range('E21').Value = '  -0.41%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '5.445'
$cell.Style = "Normal"
$ws.Range('E22').Value = '  -6.46%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '27.288.58'
$cell.Style = "Normal"
$ws.Range('E23').Value = '  -4.10%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '10.74'
$cell.Style = "Normal"
$ws.Range('E24').Value = '  -6.89%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '2.250'
$cell.Style = "Normal"
$ws.Range('E25').Value = '  -1.78%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '2.072.66'
$cell.Style = "Normal"
$ws.Range('E26').Value = '  -4.76%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '151.29'
$cell.Style = "Normal"
$ws.Range('E27').Value = '  -2.50%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '19.31'
$cell.Style = "Normal"
$ws.Range('E28').Value = '  -4.33%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '2.043'
$cell.Style = "Normal"
$ws.Range('E29').Value = '  -5.15%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '5.501'
$cell.Style = "Normal"
$ws.Range('E30').Value = '  -6.98%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '119.02'
$cell.Style = "Normal"
$ws.Range('E31').Value = '  -3.64%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '0.09291'
$cell.Style = "Normal"
$ws.Range('E32').Value = '  -3.15%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '1.451'
$cell.Style = "Normal"
$ws.Range('E33').Value = '  +0.20%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '0.9243'
$cell.Style = "Normal"
$ws.Range('E34').Value = '  -5.52%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '3.622'
$cell.Style = "Normal"
$ws.Range('E35').Value = '  -0.99%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '5.230'
$cell.Style = "Normal"
$ws.Range('E36').Value = '  -6.59%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '0.02205'
$cell.Style = "Normal"
$ws.Range('E37').Value = '  -5.47%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '0.05932'
$cell.Style = "Normal"
$ws.Range('E38').Value = '  -4.45%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '1.198'
$cell.Style = "Normal"
$ws.Range('E39').Value = '  -3.53%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '8.250'
$cell.Style = "Normal"
$ws.Range('E40').Value = '  -7.36%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '1.003'
$cell.Style = "Normal"
$ws.Range('E41').Value = '  -0.37%  '
$ws.Range('E42').Value = '  -5.39%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '0.1844'
$cell.Style = "Normal"
$ws.Range('E43').Value = '  -4.00%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '10.23'
$cell.Style = "Normal"
$ws.Range('E44').Value = '  -8.11%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '1.253'
$cell.Style = "Normal"
$ws.Range('E45').Value = '  -7.72%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '0.5609'
$cell.Style = "Normal"
$ws.Range('E46').Value = '  -5.78%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '12.02'
$cell.Style = "Normal"
$ws.Range('E47').Value = '  -7.85%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '3.355'
$cell.Style = "Normal"
$ws.Range('E48').Value = '  -1.12%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '1.910'
$cell.Style = "Normal"
$ws.Range('E49').Value = '  -7.20%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '0.06847'
$cell.Style = "Normal"
$ws.Range('E50').Value = '  +0.31%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '108.10'
$cell.Style = "Normal"
$ws.Range('E51').Value = '  -2.00%  '
